$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2172774869109948
$ws.Range("C2").Value = 0.5209424083769634
$ws.Range("J2").Value = 0.01308900523560209
$ws.Range("P2").Value = 0.1465968586387434
$ws.Range("S2").Value = 0.1020942408376963
$ws.Range("B3").Value = 0.009523809523809525
$ws.Range("C3").Value = 0.02380952380952381
$ws.Range("J3").Value = 0.07142857142857142
$ws.Range("P3").Value = 0.719047619047619
$ws.Range("S3").Value = 0.1761904761904762
$ws.Range("J4").Value = 0.02
$ws.Range("P4").Value = 0.74
$ws.Range("S4").Value = 0.24
$ws.Range("J5").Value = 0.5
$ws.Range("P5").Value = 0.5
$ws.Range("B6").Value = 0.06018518518518518
$ws.Range("D6").Value = 0.009259259259259259
$ws.Range("E6").Value = 0.004629629629629629
$ws.Range("F6").Value = 0.05092592592592592
$ws.Range("J6").Value = 0.2592592592592592
$ws.Range("O6").Value = 0.03240740740740741
$ws.Range("Q6").Value = 0.1527777777777778
$ws.Range("R6").Value = 0.07407407407407407
$ws.Range("S6").Value = 0.3564814814814815
$ws.Range("B7").Value = 0.1167315175097276
$ws.Range("D7").Value = 0.01945525291828794
$ws.Range("F7").Value = 0.05836575875486381
$ws.Range("J7").Value = 0.1673151750972763
$ws.Range("O7").Value = 0.01945525291828794
$ws.Range("Q7").Value = 0.1867704280155642
$ws.Range("R7").Value = 0.07003891050583658
$ws.Range("S7").Value = 0.3618677042801556
$ws.Range("B8").Value = 0.09183673469387756
$ws.Range("D8").Value = 0.02551020408163265
$ws.Range("F8").Value = 0.03316326530612245
$ws.Range("J8").Value = 0.1581632653061225
$ws.Range("O8").Value = 0.01275510204081633
$ws.Range("Q8").Value = 0.2219387755102041
$ws.Range("R8").Value = 0.06377551020408163
$ws.Range("S8").Value = 0.3928571428571428
$ws.Range("B9").Value = 0.1265060240963855
$ws.Range("F9").Value = 0.0783132530120482
$ws.Range("J9").Value = 0.0963855421686747
$ws.Range("O9").Value = 0.03012048192771084
$ws.Range("Q9").Value = 0.2108433734939759
$ws.Range("R9").Value = 0.04216867469879518
$ws.Range("S9").Value = 0.4156626506024096
$ws.Range("B10").Value = 0.1269230769230769
$ws.Range("D10").Value = 0.02243589743589744
$ws.Range("E10").Value = 0.000641025641025641
$ws.Range("F10").Value = 0.06602564102564103
$ws.Range("J10").Value = 0.1403846153846154
$ws.Range("O10").Value = 0.01153846153846154
$ws.Range("Q10").Value = 0.2307692307692308
$ws.Range("R10").Value = 0.05961538461538462
$ws.Range("S10").Value = 0.3416666666666667
$ws.Range("G11").Value = 0.1418269230769231
$ws.Range("J11").Value = 0.09134615384615384
$ws.Range("K11").Value = 0.1826923076923077
$ws.Range("L11").Value = 0.5673076923076923
$ws.Range("S11").Value = 0.01682692307692308
$ws.Range("G12").Value = 0.6916666666666667
$ws.Range("J12").Value = 0.2708333333333333
$ws.Range("K12").Value = 0.01666666666666667
$ws.Range("L12").Value = 0.008333333333333333
$ws.Range("S12").Value = 0.0125
$ws.Range("F13").Value = 0.01886792452830189
$ws.Range("G13").Value = 0.7358490566037735
$ws.Range("J13").Value = 0.2264150943396226
$ws.Range("S13").Value = 0.01886792452830189
$ws.Range("F15").Value = 0.01214574898785425
$ws.Range("H15").Value = 0.1336032388663968
$ws.Range("I15").Value = 0.06882591093117409
$ws.Range("J15").Value = 0.4412955465587045
$ws.Range("K15").Value = 0.06072874493927125
$ws.Range("M15").Value = 0.02024291497975709
$ws.Range("O15").Value = 0.048582995951417
$ws.Range("S15").Value = 0.2145748987854251
$ws.Range("F16").Value = 0.01694915254237288
$ws.Range("H16").Value = 0.1440677966101695
$ws.Range("I16").Value = 0.05932203389830509
$ws.Range("J16").Value = 0.4915254237288136
$ws.Range("K16").Value = 0.1355932203389831
$ws.Range("M16").Value = 0.0211864406779661
$ws.Range("O16").Value = 0.05084745762711865
$ws.Range("S16").Value = 0.08050847457627118
$ws.Range("F17").Value = 0.01773049645390071
$ws.Range("H17").Value = 0.1436170212765958
$ws.Range("I17").Value = 0.07624113475177305
$ws.Range("J17").Value = 0.4556737588652482
$ws.Range("K17").Value = 0.1294326241134752
$ws.Range("M17").Value = 0.01773049645390071
$ws.Range("N17").Value = 0.003546099290780142
$ws.Range("O17").Value = 0.06560283687943262
$ws.Range("S17").Value = 0.09042553191489362
$ws.Range("F18").Value = 0.006329113924050633
$ws.Range("H18").Value = 0.1518987341772152
$ws.Range("I18").Value = 0.05063291139240506
$ws.Range("J18").Value = 0.4430379746835443
$ws.Range("K18").Value = 0.1582278481012658
$ws.Range("M18").Value = 0.03164556962025317
$ws.Range("O18").Value = 0.06329113924050633
$ws.Range("S18").Value = 0.0949367088607595
$ws.Range("F19").Value = 0.01410658307210031
$ws.Range("H19").Value = 0.1724137931034483
$ws.Range("I19").Value = 0.06661442006269593
$ws.Range("J19").Value = 0.3894984326018809
$ws.Range("K19").Value = 0.1473354231974922
$ws.Range("M19").Value = 0.02351097178683386
$ws.Range("O19").Value = 0.08307210031347963
$ws.Range("S19").Value = 0.103448275862069
